$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 91) with the latest Adafruit IO reading,
# following the same layout as the existing rows:
# Timestamp | Feed Key | Value | Latitude | Longitude | Elevation
$newRow = 91

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"
$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"

# Column C holds a numeric-looking reading ("25") that must still be
# stored as plain text, matching the rest of the sheet's data rows.
# Writing it directly as a string gets auto-coerced to a number, so
# instead compute it as text via a formula and then paste the result
# back as a literal value (no formula, no extra cell style left behind).
$c = $ws.Cells.Item($newRow, 3)
$c.Formula = '=TEXT(25,"0")'
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
